$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 302
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45172
